# "Criando Forms web rpa" - refresh the CNPJ lookup sheet with a real list of
# CNPJ numbers to process, and trim the "Dados" (results) sheet back down to
# just its header row, ready to be (re)filled by the RPA run.

$wb = $excel.ActiveWorkbook
$wsDados = $wb.Worksheets.Item("Dados")
$wsCnpjs = $wb.Worksheets.Item("CNPJS")

# CNPJS sheet: replace the sample CNPJ entries with the batch to look up.
# Row 1 (header "CNPJ") is left untouched.
$wsCnpjs.Range("A2").Value = "14.380.200/0001-21"
$wsCnpjs.Range("A3").Value = "43.828.151/0001-45"
$wsCnpjs.Range("A4").Value = "29.979.036/0001-40"
$wsCnpjs.Range("A5").Value = "60.316.817/0001-03"

# Dados sheet: drop the sample result row, keeping only the header row
# (CNPJ, nome, situacao, atividade_principal, cep, email).
$wsDados.Rows("2").Delete()

# Widen the "nome" and "atividade_principal" columns so the long values the
# RPA run will populate (company name / business activity description) are
# readable without truncation.
$wsDados.Columns.Item(2).ColumnWidth = 61.6002380952381
$wsDados.Columns.Item(4).ColumnWidth = 21.028809523809525
